$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 35: 2014-02-22 (serial 41692), 0.25 hours ---
$ws.Range("A34").Copy()
$ws.Range("A35").PasteSpecial(-4122)  # xlPasteFormats, keep existing numeric/date formatting
$ws.Range("A35").Value2 = 41692
$ws.Range("B35").Value = "Updated the software requirements specification document according to the team's criticisms."
$ws.Range("I35").Value = 0.25

# --- Row 36: 2014-02-25 (serial 41695), 1.5 hours ---
$ws.Range("A34").Copy()
$ws.Range("A36").PasteSpecial(-4122)
$ws.Range("A36").Value2 = 41695
$ws.Range("B36").Value = "Client meeting. Recorded meeting minutes and client feedback. Also recorded team milestones."
$ws.Range("I36").Value = 1.5

# --- Row 37: 2014-02-27 (serial 41697), 1 hour ---
$ws.Range("A34").Copy()
$ws.Range("A37").PasteSpecial(-4122)
$ws.Range("A37").Value2 = 41697
$ws.Range("B37").Value = "Finished recording meeting minutes. Updated the software requirements specification document in accordance with the client's feedback."
$ws.Range("I37").Value = 1

$excel.CutCopyMode = 0

# Update selection to reflect the new active cell/range used while editing the log
$ws.Range("B38:H38").Select()
